$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (C1:E1) -------------------------------------------
$ws.Range("C1").Value = "Processing Time (s)"
$ws.Range("D1").Value = "Processing Time per File (s)"
$ws.Range("E1").Value = "Cost per Million Tokens"

# Match the existing header formatting (bold / centered / bordered), the
# same way Excel users do it interactively: copy an already-styled header
# cell and paste only its formats onto the new header cells.
$ws.Range("B1").Copy()
$ws.Range("C1:E1").PasteSpecial(-4122)

# --- Row 2 (Presidio) new metric columns --------------------------------
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

# --- Row 3 (llama-3-70b): updated PII value + new metric columns -------
$ws.Range("B3").Value = 78
$ws.Range("C3").Value = 733
$ws.Range("D3").Value = 7.33
$ws.Range("E3").Value = 0.81

# --- Row 4 (llama-3-8b): updated PII value + new metric columns --------
$ws.Range("B4").Value = 74
$ws.Range("C4").Value = 406
$ws.Range("D4").Value = 4.06
$ws.Range("E4").Value = 0.18

# --- Row 5 (qwen-110b): updated PII value + new metric columns ---------
$ws.Range("B5").Value = 69
$ws.Range("C5").Value = 639
$ws.Range("D5").Value = 6.39
$ws.Range("E5").Value = 1.62
